$wb = $excel.ActiveWorkbook

# --- Budget Out sheet ---
$wsBudgetOut = $wb.Worksheets.Item("Budget Out")
$wsBudgetOut.Range("C9").Value = 91.42
$wsBudgetOut.Range("F9").Value = "Description007zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# --- TestRecord sheet ---
$wsTestRecord = $wb.Worksheets.Item("TestRecord")
$wsTestRecord.Range("A10").Value = 43264
$wsTestRecord.Range("B10").Value = 122.34
$wsTestRecord.Range("E10").Value = "some test textzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# --- Expected Out sheet ---
$wsExpectedOut = $wb.Worksheets.Item("Expected Out")
$wsExpectedOut.Range("B9").Value = 1349.36
$wsExpectedOut.Range("B11").Value = 429.22
# B1 holds =SUM(B2:B295); changing B9/B11 (which are inside that range)
# causes Excel to recalc B1 automatically to the new total.
